# Auto-generated edit script applying the Spriggan_Profits (per-class Leve profit sheets) data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 30
$ws.Range("I4").Value = 30
$ws.Range("K4").Value = 30
$ws.Range("M4").Value = 84

$ws.Range("H39").Value = 5639.5557
$ws.Range("I39").Value = 64.25
$ws.Range("J39").Value = 10099.8
$ws.Range("K39").Value = 192.75
$ws.Range("L39").Value = 30299.4
$ws.Range("M39").Value = 103.25
$ws.Range("N39").Value = -30891.4

$ws.Range("H62").Value = 4902.3335
$ws.Range("J62").Value = 4710
$ws.Range("L62").Value = 4710
$ws.Range("N62").Value = -5958

$ws.Range("H64").Value = 20837428
$ws.Range("I64").Value = 33336934
$ws.Range("K64").Value = 33336934
$ws.Range("M64").Value = -33336686

$ws.Range("H65").Value = 4902.3335
$ws.Range("J65").Value = 4710
$ws.Range("L65").Value = 23550
$ws.Range("N65").Value = -29790

$ws.Range("H67").Value = 20837428
$ws.Range("I67").Value = 33336934
$ws.Range("K67").Value = 33336934
$ws.Range("M67").Value = -33336076

$ws.Range("H112").Value = 70206.266
$ws.Range("I112").Value = 1750
$ws.Range("J112").Value = 75096
$ws.Range("K112").Value = 5250
$ws.Range("L112").Value = 225288
$ws.Range("M112").Value = -4142
$ws.Range("N112").Value = -227504

$ws.Range("H132").Value = 2572.3513
$ws.Range("I132").Value = 2629.3235
$ws.Range("K132").Value = 7887.970499999999
$ws.Range("M132").Value = -5357.970499999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11096.352
$ws.Range("I32").Value = 7641.1816
$ws.Range("J32").Value = 26299.1
$ws.Range("K32").Value = 7641.1816
$ws.Range("L32").Value = 26299.1
$ws.Range("M32").Value = -7354.1816
$ws.Range("N32").Value = -26873.1

$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()

$ws.Range("H74").Value = 26320286
$ws.Range("I74").Value = 28576106
$ws.Range("J74").Value = 2384.3333
$ws.Range("K74").Value = 28576106
$ws.Range("L74").Value = 2384.3333
$ws.Range("M74").Value = -28575232
$ws.Range("N74").Value = -4132.3333

$ws.Range("H77").Value = 26320286
$ws.Range("I77").Value = 28576106
$ws.Range("J77").Value = 2384.3333
$ws.Range("K77").Value = 142880530
$ws.Range("L77").Value = 11921.6665
$ws.Range("M77").Value = -142876162
$ws.Range("N77").Value = -20657.6665

$ws.Range("H110").Value = 64879
$ws.Range("I110").Value = 73360.78999999999
$ws.Range("J110").Value = 5506.5
$ws.Range("K110").Value = 73360.78999999999
$ws.Range("L110").Value = 5506.5
$ws.Range("M110").Value = -71315.78999999999
$ws.Range("N110").Value = -9596.5

$ws.Range("H122").Value = 2078.0356
$ws.Range("I122").Value = 2225.913
$ws.Range("K122").Value = 6677.739
$ws.Range("M122").Value = -4227.739

$ws.Range("H132").Value = 3452398.8
$ws.Range("I132").Value = 4170006.2
$ws.Range("J132").Value = 7882.6
$ws.Range("K132").Value = 12510018.6
$ws.Range("L132").Value = 23647.8
$ws.Range("M132").Value = -12507488.6
$ws.Range("N132").Value = -28707.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1387.7812
$ws.Range("I20").Value = 1319.9546
$ws.Range("J20").Value = 1537
$ws.Range("K20").Value = 1319.9546
$ws.Range("L20").Value = 1537
$ws.Range("M20").Value = -1072.9546
$ws.Range("N20").Value = -2031

$ws.Range("H107").Value = 32301.875
$ws.Range("I107").Value = 1190.1428
$ws.Range("K107").Value = 1190.1428
$ws.Range("M107").Value = 729.8571999999999

$ws.Range("H134").Value = 12197866
$ws.Range("I134").Value = 13516042
$ws.Range("J134").Value = 4743
$ws.Range("K134").Value = 40548126
$ws.Range("L134").Value = 14229
$ws.Range("M134").Value = -40545591
$ws.Range("N134").Value = -19299

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5420.528
$ws.Range("J31").Value = 8076.077
$ws.Range("L31").Value = 8076.077
$ws.Range("N31").Value = -8666.077000000001

$ws.Range("H34").Value = 5420.528
$ws.Range("J34").Value = 8076.077
$ws.Range("L34").Value = 8076.077
$ws.Range("N34").Value = -8480.077000000001

$ws.Range("H43").Value = 23915
$ws.Range("J43").Value = 26898
$ws.Range("L43").Value = 26898
$ws.Range("N43").Value = -27266

$ws.Range("H69").Value = 3660.6667
$ws.Range("I69").Value = 3660.6667
$ws.Range("K69").Value = 3660.6667
$ws.Range("M69").Value = -2911.6667

$ws.Range("H72").Value = 3660.6667
$ws.Range("I72").Value = 3660.6667
$ws.Range("K72").Value = 10982.0001
$ws.Range("M72").Value = -7238.000100000001

$ws.Range("H101").Value = 23915
$ws.Range("J101").Value = 26898
$ws.Range("L101").Value = 26898
$ws.Range("N101").Value = -33388

$ws.Range("H132").Value = 24391968
$ws.Range("I132").Value = 26317492
$ws.Range("K132").Value = 78952476
$ws.Range("M132").Value = -78949946

$ws.Range("H141").Value = 202798
$ws.Range("J141").Value = 293906.8
$ws.Range("L141").Value = 293906.8
$ws.Range("N141").Value = -304266.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 457.14285
$ws.Range("I46").Value = 460
$ws.Range("K46").Value = 1380
$ws.Range("M46").Value = -1289

$ws.Range("H114").Value = 73933.86
$ws.Range("I114").Value = 167158.83
$ws.Range("J114").Value = 4015.125
$ws.Range("K114").Value = 501476.49
$ws.Range("L114").Value = 12045.375
$ws.Range("M114").Value = -498222.49
$ws.Range("N114").Value = -18553.375

$ws.Range("H131").Value = 2394.5
$ws.Range("J131").Value = 2493.125
$ws.Range("L131").Value = 7479.375
$ws.Range("N131").Value = -17559.375

$ws.Range("H137").Value = 12504335
$ws.Range("J137").Value = 3497.5
$ws.Range("L137").Value = 10492.5
$ws.Range("N137").Value = -20692.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2994.4517
$ws.Range("J80").Value = 2865.2693
$ws.Range("L80").Value = 2865.2693
$ws.Range("N80").Value = -4861.2693

$ws.Range("H83").Value = 2994.4517
$ws.Range("J83").Value = 2865.2693
$ws.Range("L83").Value = 14326.3465
$ws.Range("N83").Value = -24310.3465

$ws.Range("H97").Value = 1043.1
$ws.Range("I97").Value = 805.8333
$ws.Range("J97").Value = 1399
$ws.Range("K97").Value = 805.8333
$ws.Range("L97").Value = 1399
$ws.Range("M97").Value = -309.8333
$ws.Range("N97").Value = -2391

$ws.Range("H99").Value = 18461.3
$ws.Range("I99").Value = 8076.75
$ws.Range("K99").Value = 8076.75
$ws.Range("M99").Value = -5830.75

$ws.Range("H102").Value = 2467.7646
$ws.Range("I102").Value = 2372
$ws.Range("J102").Value = 4000
$ws.Range("K102").Value = 2372
$ws.Range("L102").Value = 4000
$ws.Range("M102").Value = -750
$ws.Range("N102").Value = -7244

$ws.Range("H107").Value = 2898.077
$ws.Range("I107").Value = 2130.6667
$ws.Range("K107").Value = 2130.6667
$ws.Range("M107").Value = -210.6667000000002

$ws.Range("H126").Value = 7007.077
$ws.Range("I126").Value = 7452.5713
$ws.Range("J126").Value = 6487.3335
$ws.Range("K126").Value = 22357.7139
$ws.Range("L126").Value = 19462.0005
$ws.Range("M126").Value = -19887.7139
$ws.Range("N126").Value = -24402.0005

$ws.Range("H132").Value = 3293251
$ws.Range("I132").Value = 3382021.5
$ws.Range("K132").Value = 10146064.5
$ws.Range("M132").Value = -10143534.5

$ws.Range("H141").Value = 84990
$ws.Range("J141").Value = 84990
$ws.Range("L141").Value = 84990
$ws.Range("N141").Value = -95350

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 587.1177
$ws.Range("J55").Value = 844.1667
$ws.Range("L55").Value = 844.1667
$ws.Range("N55").Value = -1190.1667

$ws.Range("H99").Value = 40313.668
$ws.Range("I99").Value = 25899.5
$ws.Range("K99").Value = 25899.5
$ws.Range("M99").Value = -22904.5

$ws.Range("H100").Value = 9014306
$ws.Range("I100").Value = 10437312
$ws.Range("J100").Value = 1932.6666
$ws.Range("K100").Value = 10437312
$ws.Range("L100").Value = 1932.6666
$ws.Range("M100").Value = -10436771
$ws.Range("N100").Value = -3014.6666

$ws.Range("H132").Value = 39453310
$ws.Range("I132").Value = 78902740
$ws.Range("J132").Value = 3888.889
$ws.Range("K132").Value = 236708220
$ws.Range("L132").Value = 11666.667
$ws.Range("M132").Value = -236705690
$ws.Range("N132").Value = -16726.667
